$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 - copy of row 5 with updated date
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44461
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100107
$ws.Cells.Item(8, 8).Value = "Otros"
$ws.Cells.Item(8, 9).Value = 100107002
$ws.Cells.Item(8, 10).Value = "Chirimoya"
$ws.Cells.Item(8, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(8, 12).Value = "Especial"
$ws.Cells.Item(8, 13).Value = 60
$ws.Cells.Item(8, 14).Value = 31000
$ws.Cells.Item(8, 15).Value = 32000
$ws.Cells.Item(8, 16).Value = 31500
$ws.Cells.Item(8, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value = 3150
$ws.Cells.Item(8, 20).Value = 10

# New row 9 - copy of row 6 with updated date
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44461
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100107
$ws.Cells.Item(9, 8).Value = "Otros"
$ws.Cells.Item(9, 9).Value = 100107002
$ws.Cells.Item(9, 10).Value = "Chirimoya"
$ws.Cells.Item(9, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 30
$ws.Cells.Item(9, 14).Value = 30000
$ws.Cells.Item(9, 15).Value = 30000
$ws.Cells.Item(9, 16).Value = 30000
$ws.Cells.Item(9, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 19).Value = 3000
$ws.Cells.Item(9, 20).Value = 10
